$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 563.64105
$ws.Range("J17").Value = 563.64105
$ws.Range("L17").Value = 1690.92315
$ws.Range("N17").Value = -2026.92315
$ws.Range("H43").Value = 1678.32
$ws.Range("I43").Value = 903.1111
$ws.Range("J43").Value = 2114.375
$ws.Range("K43").Value = 903.1111
$ws.Range("L43").Value = 2114.375
$ws.Range("M43").Value = -834.1111
$ws.Range("N43").Value = -2252.375
$ws.Range("H70").Value = 2834.2222
$ws.Range("I70").Value = 3246.25
$ws.Range("K70").Value = 9738.75
$ws.Range("M70").Value = -9468.75
$ws.Range("H73").Value = 2834.2222
$ws.Range("I73").Value = 3246.25
$ws.Range("K73").Value = 9738.75
$ws.Range("M73").Value = -8802.75
$ws.Range("H97").Value = 2000
$ws.Range("J97").Value = 2000
$ws.Range("L97").Value = 6000
$ws.Range("N97").Value = -6992
$ws.Range("H112").Value = 2711792.5
$ws.Range("J112").Value = 2925811.5
$ws.Range("L112").Value = 8777434.5
$ws.Range("N112").Value = -8779650.5
$ws.Range("H132").Value = 1066623.4
$ws.Range("I132").Value = 1260.8718
$ws.Range("K132").Value = 3782.6154
$ws.Range("M132").Value = -1252.6154
$ws.Range("H135").Value = 17171.887
$ws.Range("I135").Value = 21849.299
$ws.Range("J135").Value = 2516
$ws.Range("K135").Value = 196643.691
$ws.Range("L135").Value = 22644
$ws.Range("M135").Value = -194108.691
$ws.Range("N135").Value = -27714
$ws.Range("H141").Value = 656.1429000000001
$ws.Range("I141").Value = 579.8387
$ws.Range("J141").Value = 1247.5
$ws.Range("K141").Value = 1739.5161
$ws.Range("L141").Value = 3742.5
$ws.Range("M141").Value = 3440.4839
$ws.Range("N141").Value = -14102.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3473.2
$ws.Range("I2").Value = 3672
$ws.Range("J2").Value = 3175
$ws.Range("K2").Value = 3672
$ws.Range("L2").Value = 3175
$ws.Range("M2").Value = -3559
$ws.Range("N2").Value = -3401
$ws.Range("H14").Value = 2851.75
$ws.Range("I14").Value = 3400
$ws.Range("J14").Value = 1207
$ws.Range("K14").Value = 3400
$ws.Range("L14").Value = 1207
$ws.Range("M14").Value = -3225
$ws.Range("N14").Value = -1557
$ws.Range("H32").Value = 1360.21
$ws.Range("I32").Value = 1081.0952
$ws.Range("J32").Value = 2825.5625
$ws.Range("K32").Value = 1081.0952
$ws.Range("L32").Value = 2825.5625
$ws.Range("M32").Value = -794.0952
$ws.Range("N32").Value = -3399.5625
$ws.Range("H61").Value = 22773450
$ws.Range("I61").Value = 26342770
$ws.Range("J61").Value = 167764.17
$ws.Range("K61").Value = 26342770
$ws.Range("L61").Value = 167764.17
$ws.Range("M61").Value = -26342558
$ws.Range("N61").Value = -168188.17
$ws.Range("H74").Value = 4943041.5
$ws.Range("I74").Value = 6123273.5
$ws.Range("J74").Value = 104091.2
$ws.Range("K74").Value = 6123273.5
$ws.Range("L74").Value = 104091.2
$ws.Range("M74").Value = -6122399.5
$ws.Range("N74").Value = -105839.2
$ws.Range("H77").Value = 4943041.5
$ws.Range("I77").Value = 6123273.5
$ws.Range("J77").Value = 104091.2
$ws.Range("K77").Value = 30616367.5
$ws.Range("L77").Value = 520456
$ws.Range("M77").Value = -30611999.5
$ws.Range("N77").Value = -529192
$ws.Range("H116").Value = 3473.2
$ws.Range("I116").Value = 3672
$ws.Range("J116").Value = 3175
$ws.Range("K116").Value = 3672
$ws.Range("L116").Value = 3175
$ws.Range("M116").Value = -1378
$ws.Range("N116").Value = -7763
$ws.Range("H136").Value = 22773450
$ws.Range("I136").Value = 26342770
$ws.Range("J136").Value = 167764.17
$ws.Range("K136").Value = 79028310
$ws.Range("L136").Value = 503292.51
$ws.Range("M136").Value = -79025760
$ws.Range("N136").Value = -508392.51

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3473.2
$ws.Range("I3").Value = 3672
$ws.Range("J3").Value = 3175
$ws.Range("K3").Value = 3672
$ws.Range("L3").Value = 3175
$ws.Range("M3").Value = -3558
$ws.Range("N3").Value = -3403

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 3000
$ws.Range("I25").Value = 3000
$ws.Range("K25").Value = 3000
$ws.Range("M25").Value = -2826
$ws.Range("H58").Value = 23810628
$ws.Range("I58").Value = 30303872
$ws.Range("J58").Value = 2068
$ws.Range("K58").Value = 30303872
$ws.Range("L58").Value = 2068
$ws.Range("M58").Value = -30303669
$ws.Range("N58").Value = -2474
$ws.Range("H64").Value = 28628.334
$ws.Range("J64").Value = 28628.334
$ws.Range("L64").Value = 28628.334
$ws.Range("N64").Value = -29124.334
$ws.Range("H67").Value = 28628.334
$ws.Range("J67").Value = 28628.334
$ws.Range("L67").Value = 28628.334
$ws.Range("N67").Value = -30344.334
$ws.Range("H134").Value = 122645.78
$ws.Range("I134").Value = 3718.6667
$ws.Range("J134").Value = 360500
$ws.Range("K134").Value = 11156.0001
$ws.Range("L134").Value = 1081500
$ws.Range("M134").Value = -8621.000100000001
$ws.Range("N134").Value = -1086570
$ws.Range("H136").Value = 23810628
$ws.Range("I136").Value = 30303872
$ws.Range("J136").Value = 2068
$ws.Range("K136").Value = 90911616
$ws.Range("L136").Value = 6204
$ws.Range("M136").Value = -90909066
$ws.Range("N136").Value = -11304

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 385154.62
$ws.Range("I11").Value = 572.8570999999999
$ws.Range("J11").Value = 833833.3
$ws.Range("K11").Value = 1718.5713
$ws.Range("L11").Value = 2501499.9
$ws.Range("M11").Value = -1578.5713
$ws.Range("N11").Value = -2501779.9
$ws.Range("H69").Value = 898.5
$ws.Range("I69").Value = 500
$ws.Range("J69").Value = 978.2
$ws.Range("K69").Value = 1500
$ws.Range("L69").Value = 2934.6
$ws.Range("M69").Value = -689
$ws.Range("N69").Value = -4556.6
$ws.Range("H72").Value = 898.5
$ws.Range("I72").Value = 500
$ws.Range("J72").Value = 978.2
$ws.Range("K72").Value = 4500
$ws.Range("L72").Value = 8803.800000000001
$ws.Range("M72").Value = -444
$ws.Range("N72").Value = -16915.8
$ws.Range("H129").Value = 5954927
$ws.Range("I129").Value = 2882
$ws.Range("J129").Value = 9261619
$ws.Range("K129").Value = 8646
$ws.Range("L129").Value = 27784857
$ws.Range("M129").Value = -3646
$ws.Range("N129").Value = -27794857
$ws.Range("H131").Value = 13514438
$ws.Range("J131").Value = 1045.5574
$ws.Range("L131").Value = 3136.6722
$ws.Range("N131").Value = -13216.6722

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H38").Value = 9800
$ws.Range("J38").Value = 9800
$ws.Range("L38").Value = 9800
$ws.Range("N38").Value = -10726
$ws.Range("H122").Value = 2899.647
$ws.Range("I122").Value = 2636.24
$ws.Range("J122").Value = 3631.3333
$ws.Range("K122").Value = 7908.719999999999
$ws.Range("L122").Value = 10893.9999
$ws.Range("M122").Value = -5458.719999999999
$ws.Range("N122").Value = -15793.9999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 996.0833
$ws.Range("J9").Value = 1523.4286
$ws.Range("L9").Value = 1523.4286
$ws.Range("N9").Value = -1971.4286
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("M13").ClearContents()
$ws.Range("H68").Value = 1698.6957
$ws.Range("I68").Value = 1653.5
$ws.Range("K68").Value = 1653.5
$ws.Range("M68").Value = -904.5
$ws.Range("H71").Value = 1698.6957
$ws.Range("I71").Value = 1653.5
$ws.Range("K71").Value = 8267.5
$ws.Range("M71").Value = -4523.5

Write-Output "Done applying changes"
